# project setup for training images
# Update the "Vision Statement Reminder" slide body copy: drop the
# "multi-platform" qualifier from the app description.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$body = $s.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "We are creating a mobile app which focuses on image recognition capabilities in an educational environment. "
